$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value, derived from the updated
# cryptocurrency price/volume snapshot (and the BabyDogeCoin rank change
# that shifted rows 47-51 down by one).
$cellValues = [ordered]@{
    "D2" = "29.262.31"
    "E2" = "  -0.61%  "
    "D3" = "1.838.97"
    "E3" = "  -0.62%  "
    "D4" = "0.9983"
    "E4" = "  -0.21%  "
    "D5" = "239.64"
    "E5" = "  -0.31%  "
    "D6" = "0.6249"
    "E6" = "  -0.45%  "
    "D7" = "0.9986"
    "E7" = "  -0.25%  "
    "D8" = "0.07474"
    "E8" = "  -2.60%  "
    "E9" = "  -0.92%  "
    "D10" = "24.22"
    "E10" = "  -2.77%  "
    "E11" = "  -0.56%  "
    "D12" = "1.839.03"
    "E12" = "  -0.80%  "
    "D13" = "4.980"
    "E13" = "  -1.12%  "
    "D14" = "0.6763"
    "E14" = "  -0.81%  "
    "D15" = "0.00001027"
    "E15" = "  -4.22%  "
    "D16" = "82.00"
    "E16" = "  -1.79%  "
    "D17" = "2.099.66"
    "E17" = "  -0.43%  "
    "D18" = "6.085"
    "E18" = "  -2.17%  "
    "D19" = "29.289.26"
    "E19" = "  -0.54%  "
    "D20" = "228.29"
    "E20" = "  -0.19%  "
    "D21" = "12.24"
    "E21" = "  -1.24%  "
    "D22" = "0.9989"
    "E22" = "  -0.22%  "
    "D23" = "7.368"
    "E23" = "  -1.03%  "
    "D24" = "0.9994"
    "E24" = "  -0.12%  "
    "D25" = "157.99"
    "E25" = "  +0.07%  "
    "D26" = "0.1377"
    "E26" = "  -0.12%  "
    "D27" = "8.360"
    "E27" = "  -0.71%  "
    "E28" = "  -1.29%  "
    "E29" = "  +2.16%  "
    "D30" = "1.469"
    "E30" = "  +0.57%  "
    "D31" = "0.05700"
    "E31" = "  +1.45%  "
    "D32" = "4.092"
    "E32" = "  -0.87%  "
    "D33" = "4.024"
    "E33" = "  -0.57%  "
    "D34" = "1.816"
    "E34" = "  -1.54%  "
    "E35" = "  -1.95%  "
    "D36" = "0.6923"
    "E36" = "  -2.07%  "
    "D37" = "2.583"
    "E37" = "  -0.56%  "
    "D38" = "2.818"
    "E38" = "  +2.25%  "
    "D39" = "1.241.30"
    "E39" = "  +1.30%  "
    "D40" = "0.01807"
    "E40" = "  +0.84%  "
    "D41" = "6.500"
    "E41" = "  +0.81%  "
    "D42" = "0.9035"
    "E42" = "  +0.10%  "
    "D43" = "0.9974"
    "E43" = "  -0.38%  "
    "D44" = "2.001.37"
    "E44" = "  -0.67%  "
    "D45" = "101.26"
    "E45" = "  -0.72%  "
    "D46" = "65.60"
    "E46" = "  -0.88%  "
    "B47" = "BabyDogeCoin"
    "C47" = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
    "D47" = "0.00000000119"
    "E47" = "  -1.20%  "
    "B48" = "Aptos"
    "C48" = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
    "D48" = "7.054"
    "E48" = "  -2.13%  "
    "B49" = "Algorand"
    "C49" = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
    "D49" = "0.1168"
    "E49" = "  +1.01%  "
    "B50" = "EnergySwap"
    "C50" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D50" = "8.993"
    "E50" = "  -0.45%  "
    "B51" = "TheSandbox"
    "C51" = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
    "D51" = "0.3933"
    "E51" = "  -2.20%  "
}

foreach ($ref in $cellValues.Keys) {
    $cell = $ws.Range($ref)
    # Force text format first so numeric-looking strings (prices like
    # "0.9983" or "29.262.31") are stored verbatim as text instead of
    # being auto-converted to numbers/dates by Excel.
    $cell.NumberFormat = "@"
    $cell.Value = $cellValues[$ref]
}
